# Insert a new price record as row 263 ("Fruta / hortaliza, semanal" update).
# This shifts the existing rows 263:319 down to 264:320 (dimension grows to A1:R320)
# and fills the freshly inserted row 263 with the new weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(263).Insert()

$ws.Range("A263").Value = 1
$ws.Range("B263").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C263").Value = "Arica y Parinacota"
$ws.Range("D263").Value = 44798
$ws.Range("E263").Value = 15
$ws.Range("F263").Value = 100114013
$ws.Range("G263").Value = "Zanahoria"
$ws.Range("H263").Value = "Sin especificar"
$ws.Range("I263").Value = "Primera"
$ws.Range("J263").Value = 80
$ws.Range("K263").Value = 28000
$ws.Range("L263").Value = 30000
$ws.Range("M263").Value = 29000
$ws.Range("N263").Value = "`$/saco 25 kilos"
$ws.Range("O263").Value = "Región de Arica y Parinacota"
$ws.Range("P263").Value = 1160
$ws.Range("Q263").Value = 25
$ws.Range("R263").Value = "Hortaliza"
